# Update header labels to include the "Maximum Rate" / "Default Rate" aliases
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Mrp / Maximum Rate"
$ws.Range("E1").Value = "Sale Price / Default Rate"

# Populate sample Mrp / Sale Price values for the services row (row 6)
$ws.Range("D6").Value = 500
$ws.Range("E6").Value = 450

# Widen the columns that now hold longer header text
$ws.Columns.Item(1).ColumnWidth = 13.196666666666665
$ws.Columns.Item(3).ColumnWidth = 16.686666666666667
$ws.Columns.Item(4).ColumnWidth = 17.366666666666667
$ws.Columns.Item(5).ColumnWidth = 20.276666666666667

# Restore the active selection used when the sheet was last edited
$ws.Range("E13").Select()
